$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added ifoCAST full series evaluation: the error-table rows shift up by one
# quarter (row N now holds what used to be row N+1's values), the N (count)
# column decrements accordingly, and a brand new final row (15) appears
# with freshly computed statistics. Row 1 headers and column A (quarter
# labels) are unchanged.

$ws.Range("B2").Value = 0.2400528213269932
$ws.Range("C2").Value = 0.3753181292658889
$ws.Range("D2").Value = 0.2348887549825461
$ws.Range("E2").Value = 0.4846532316848265
$ws.Range("F2").Value = 0.4369199514236829
$ws.Range("G2").Value = 14
$ws.Range("B3").Value = 0.273599112058131
$ws.Range("C3").Value = 0.3783303919539526
$ws.Range("D3").Value = 0.2266048251956944
$ws.Range("E3").Value = 0.4760302776039508
$ws.Range("F3").Value = 0.4054553575916278
$ws.Range("G3").Value = 13
$ws.Range("B4").Value = 0.323052878118673
$ws.Range("C4").Value = 0.3688046801716363
$ws.Range("D4").Value = 0.2639536099994781
$ws.Range("E4").Value = 0.5137641579552608
$ws.Range("F4").Value = 0.4172513277134101
$ws.Range("G4").Value = 12
$ws.Range("B5").Value = 0.3213959399964313
$ws.Range("C5").Value = 0.3527736814977633
$ws.Range("D5").Value = 0.2405711618554991
$ws.Range("E5").Value = 0.4904805417705163
$ws.Range("F5").Value = 0.3885915500499728
$ws.Range("G5").Value = 11
$ws.Range("B6").Value = 0.3040077233811113
$ws.Range("C6").Value = 0.3040077233811113
$ws.Range("D6").Value = 0.2525516835738419
$ws.Range("E6").Value = 0.5025452055027905
$ws.Range("F6").Value = 0.4218095774931774
$ws.Range("G6").Value = 10
$ws.Range("B7").Value = 0.3195075457514495
$ws.Range("C7").Value = 0.3344878761007612
$ws.Range("D7").Value = 0.2334815320987104
$ws.Range("E7").Value = 0.4831992674856931
$ws.Range("F7").Value = 0.3844749898822031
$ws.Range("G7").Value = 9
$ws.Range("B8").Value = 0.3427370110204659
$ws.Range("C8").Value = 0.3605818570417189
$ws.Range("D8").Value = 0.2023703037307001
$ws.Range("E8").Value = 0.4498558699524772
$ws.Range("F8").Value = 0.3114971130477678
$ws.Range("G8").Value = 8
$ws.Range("B9").Value = 0.2987280035122604
$ws.Range("C9").Value = 0.2987280035122604
$ws.Range("D9").Value = 0.1585806734357395
$ws.Range("E9").Value = 0.3982218896993729
$ws.Range("F9").Value = 0.2844280147574629
$ws.Range("G9").Value = 7
$ws.Range("B10").Value = 0.2911926455816474
$ws.Range("C10").Value = 0.3419929353797495
$ws.Range("D10").Value = 0.3451110422066679
$ws.Range("E10").Value = 0.587461524022355
$ws.Range("F10").Value = 0.5589109611011353
$ws.Range("G10").Value = 6
$ws.Range("B11").Value = 0.2508072428796626
$ws.Range("C11").Value = 0.4114187230689154
$ws.Range("D11").Value = 0.3395202396874553
$ws.Range("E11").Value = 0.5826836531836596
$ws.Range("F11").Value = 0.5880220729345085
$ws.Range("G11").Value = 5
